$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) date in column C for rows 2-6 was updated
# from 2023-10-08 (45207) to 2023-10-09 (45208).
foreach ($row in 2..6) {
    $ws.Cells.Item($row, 3).Value = 45208
}
